# Align Sector abbreviations with the Baseline workbook's full names.
# Every yearly worksheet (2000..2100) carries the same header row 3 with
# the element/material abbreviations in D3:G3 — update them on each sheet
# so the shared text is renamed everywhere.
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $cell = $ws.Range("D3")
    if ($cell.Value2 -eq "Nd") { $cell.Value = "Neodymium" }

    $cell = $ws.Range("E3")
    if ($cell.Value2 -eq "Dy") { $cell.Value = "Dysprosium" }

    $cell = $ws.Range("F3")
    if ($cell.Value2 -eq "Cu") { $cell.Value = "Copper ores and concentrates" }

    $cell = $ws.Range("G3")
    if ($cell.Value2 -eq "Si") { $cell.Value = "Raw silicon" }
}
